$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Friday of week 31 (row 20) is now done too: bump the "days done" count
# for that week and give the Fri cell (H20) the same "done" fill/border
# formatting already used by the other completed weekday cells in that
# row (copy the format from the neighbouring done cell, G20).
$ws.Range("G20").Copy()
$ws.Range("H20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I20").Value = 5

# Move the active selection back to B2, mirroring the saved view state.
$ws.Range("B2").Select()
